$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-55 (A..G):
#   A = index, B = dataset, C = method, D = type, E = hue, F = period, G = eval
$data = @(
    @(0, "mnist", "DVI", "Train", "DVI(Train)", 0, 0.97),
    @(1, "mnist", "DVI", "Test", "DVI(Test)", 0, 0.9684999999999999),
    @(2, "mnist", "DVI", "Train", "DVI(Train)", 1, 0.997),
    @(3, "mnist", "DVI", "Test", "DVI(Test)", 1, 0.995),
    @(4, "mnist", "DVI", "Train", "DVI(Train)", 2, 1),
    @(5, "mnist", "DVI", "Test", "DVI(Test)", 2, 0.998),
    @(6, "mnist", "TimeVis", "Train", "TimeVis(Train)", 0, 0.9644999999999999),
    @(7, "mnist", "TimeVis", "Test", "TimeVis(Test)", 0, 0.9615),
    @(8, "mnist", "TimeVis", "Train", "TimeVis(Train)", 1, 0.9944999999999999),
    @(9, "mnist", "TimeVis", "Test", "TimeVis(Test)", 1, 0.991),
    @(10, "mnist", "TimeVis", "Train", "TimeVis(Train)", 2, 1),
    @(11, "mnist", "TimeVis", "Test", "TimeVis(Test)", 2, 0.997),
    @(12, "mnist", "DeepDebugger", "Train", "DeepDebugger(Train)", 0, 0.9575),
    @(13, "mnist", "DeepDebugger", "Test", "DeepDebugger(Test)", 0, 0.955),
    @(14, "mnist", "DeepDebugger", "Train", "DeepDebugger(Train)", 1, 0.996),
    @(15, "mnist", "DeepDebugger", "Test", "DeepDebugger(Test)", 1, 0.993),
    @(16, "mnist", "DeepDebugger", "Train", "DeepDebugger(Train)", 2, 1),
    @(17, "mnist", "DeepDebugger", "Test", "DeepDebugger(Test)", 2, 0.9975000000000001),
    @(18, "fmnist", "DVI", "Train", "DVI(Train)", 0, 0.8605),
    @(19, "fmnist", "DVI", "Test", "DVI(Test)", 0, 0.85),
    @(20, "fmnist", "DVI", "Train", "DVI(Train)", 1, 0.9245000000000001),
    @(21, "fmnist", "DVI", "Test", "DVI(Test)", 1, 0.9195),
    @(22, "fmnist", "DVI", "Train", "DVI(Train)", 2, 0.9550000000000001),
    @(23, "fmnist", "DVI", "Test", "DVI(Test)", 2, 0.9455),
    @(24, "fmnist", "TimeVis", "Train", "TimeVis(Train)", 0, 0.841),
    @(25, "fmnist", "TimeVis", "Test", "TimeVis(Test)", 0, 0.8315),
    @(26, "fmnist", "TimeVis", "Train", "TimeVis(Train)", 1, 0.9624999999999999),
    @(27, "fmnist", "TimeVis", "Test", "TimeVis(Test)", 1, 0.953),
    @(28, "fmnist", "TimeVis", "Train", "TimeVis(Train)", 2, 0.995),
    @(29, "fmnist", "TimeVis", "Test", "TimeVis(Test)", 2, 0.978),
    @(30, "fmnist", "DeepDebugger", "Train", "DeepDebugger(Train)", 0, 0.8765000000000001),
    @(31, "fmnist", "DeepDebugger", "Test", "DeepDebugger(Test)", 0, 0.8714999999999999),
    @(32, "fmnist", "DeepDebugger", "Train", "DeepDebugger(Train)", 1, 0.9795),
    @(33, "fmnist", "DeepDebugger", "Test", "DeepDebugger(Test)", 1, 0.972),
    @(34, "fmnist", "DeepDebugger", "Train", "DeepDebugger(Train)", 2, 0.994),
    @(35, "fmnist", "DeepDebugger", "Test", "DeepDebugger(Test)", 2, 0.981),
    @(36, "cifar10", "DVI", "Train", "DVI(Train)", 0, 0.3735),
    @(37, "cifar10", "DVI", "Test", "DVI(Test)", 0, 0.369),
    @(38, "cifar10", "DVI", "Train", "DVI(Train)", 1, 0.8815),
    @(39, "cifar10", "DVI", "Test", "DVI(Test)", 1, 0.8785000000000001),
    @(40, "cifar10", "DVI", "Train", "DVI(Train)", 2, 0.966),
    @(41, "cifar10", "DVI", "Test", "DVI(Test)", 2, 0.9445),
    @(42, "cifar10", "TimeVis", "Train", "TimeVis(Train)", 0, 0.5044999999999999),
    @(43, "cifar10", "TimeVis", "Test", "TimeVis(Test)", 0, 0.5095000000000001),
    @(44, "cifar10", "TimeVis", "Train", "TimeVis(Train)", 1, 0.9370000000000001),
    @(45, "cifar10", "TimeVis", "Test", "TimeVis(Test)", 1, 0.931),
    @(46, "cifar10", "TimeVis", "Train", "TimeVis(Train)", 2, 0.9904999999999999),
    @(47, "cifar10", "TimeVis", "Test", "TimeVis(Test)", 2, 0.9775),
    @(48, "cifar10", "DeepDebugger", "Train", "DeepDebugger(Train)", 0, 0.625),
    @(49, "cifar10", "DeepDebugger", "Test", "DeepDebugger(Test)", 0, 0.6274999999999999),
    @(50, "cifar10", "DeepDebugger", "Train", "DeepDebugger(Train)", 1, 0.946),
    @(51, "cifar10", "DeepDebugger", "Test", "DeepDebugger(Test)", 1, 0.9424999999999999),
    @(52, "cifar10", "DeepDebugger", "Train", "DeepDebugger(Train)", 2, 0.9915),
    @(53, "cifar10", "DeepDebugger", "Test", "DeepDebugger(Test)", 2, 0.98)

)

$rowNum = 2
foreach ($row in $data) {
    if ($rowNum -gt 37) {
        # New rows (38-55): copy formatting of column A from an existing data row first
        $ws.Range("A2").Copy()
        $ws.Range("A$rowNum").PasteSpecial(-4122)
    }

    $ws.Cells.Item($rowNum, 1).Value2 = $row[0]
    $ws.Cells.Item($rowNum, 2).Value2 = $row[1]
    $ws.Cells.Item($rowNum, 3).Value2 = $row[2]
    $ws.Cells.Item($rowNum, 4).Value2 = $row[3]
    $ws.Cells.Item($rowNum, 5).Value2 = $row[4]
    $ws.Cells.Item($rowNum, 6).Value2 = $row[5]
    $ws.Cells.Item($rowNum, 7).Value2 = $row[6]

    $rowNum++
}

$excel.CutCopyMode = 0

Write-Output "Updated rows 2-55 successfully"
